$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table originally had 3 data rows (Original, Normalized, Standardized).
# It now needs 2 rows per feature-set (one for each set of Initial_Weights),
# so insert one extra blank row after the "Original" row and one after the
# (now shifted) "Normalized" row, giving 6 data rows total (rows 2-7).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()

$ws.Range("A2").Value = "Original"
$ws.Range("B2").Value = "MSE"
$ws.Range("C2").Value = "[0.0, 0.0, 0.0]"
$ws.Range("D2").Value = "[5.008104594931483, 5.7341994552652675, 10.599637435894216]"
$ws.Range("E2").Value = 372.092
$ws.Range("F2").Value = 372.092
$ws.Range("G2").Value = 19.2897
$ws.Range("H2").Value = 15.2952
$ws.Range("I2").Value = 3
$ws.Range("A3").Value = "Original"
$ws.Range("B3").Value = "MSE"
$ws.Range("C3").Value = "[0.035358539505612, 0.07659193483788256, 0.03375499333000125]"
$ws.Range("D3").Value = "[5.008104595111034, 5.734199441858533, 10.599637668864629]"
$ws.Range("E3").Value = 372.092
$ws.Range("F3").Value = 372.092
$ws.Range("G3").Value = 19.2897
$ws.Range("H3").Value = 15.2952
$ws.Range("I3").Value = 4
$ws.Range("A4").Value = "Normalized"
$ws.Range("B4").Value = "MSE"
$ws.Range("C4").Value = "[0.0, 0.0, 0.0]"
$ws.Range("D4").Value = "[5.008098885544916, 5.7342039700585, 638.9028167934122]"
$ws.Range("E4").Value = 372.092
$ws.Range("F4").Value = 395137
$ws.Range("G4").Value = 628.599
$ws.Range("H4").Value = 628.303
$ws.Range("I4").Value = 4
$ws.Range("A5").Value = "Normalized"
$ws.Range("B5").Value = "MSE"
$ws.Range("C5").Value = "[0.035358539505612, 0.07659193483788256, 0.03375499333000125]"
$ws.Range("D5").Value = "[5.00809888579188, 5.734203973728281, 638.9028168269581]"
$ws.Range("E5").Value = 372.092
$ws.Range("F5").Value = 395137
$ws.Range("G5").Value = 628.599
$ws.Range("H5").Value = 628.303
$ws.Range("I5").Value = 4
$ws.Range("A6").Value = "Standardized"
$ws.Range("B6").Value = "MSE"
$ws.Range("C6").Value = "[0.0, 0.0, 0.0]"
$ws.Range("D6").Value = "[5.00810459502517, 5.734199448074711, 10.599637559505098]"
$ws.Range("E6").Value = 372.092
$ws.Range("F6").Value = 372.092
$ws.Range("G6").Value = 19.2897
$ws.Range("H6").Value = 15.2952
$ws.Range("I6").Value = 4
$ws.Range("A7").Value = "Standardized"
$ws.Range("B7").Value = "MSE"
$ws.Range("C7").Value = "[0.035358539505612, 0.07659193483788256, 0.03375499333000125]"
$ws.Range("D7").Value = "[5.00810459502517, 5.734199448074709, 10.59963755950514]"
$ws.Range("E7").Value = 372.092
$ws.Range("F7").Value = 372.092
$ws.Range("G7").Value = 19.2897
$ws.Range("H7").Value = 15.2952
$ws.Range("I7").Value = 5

